# "80% done with the merged scrapper bot" pressures fix:
#  - after placing the bet, the multiplier value (column C) is not scrapped
#    -> C436 (last row of the previous run) was left as the default
#       placeholder text instead of the real scrapped numeric total-bets
#       value, fix it to be numeric like every other completed row.
#  - after a successful cashout, the incremented value is not reset
#    -> the bot kept appending new scrapped rows (437-461) correctly as
#       numbers, but the very last row captured before the bug re-occurred
#       (462) again has its Total Bets value stuck as text, same symptom
#       as C436 above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the previously-broken trailing row: C436 should be numeric -------
$ws.Cells.Item(436, 3).Value = 6203

# --- Newly scrapped rows (437-462) -----------------------------------------
# Columns A (multiplier) and B (timestamp) are always scrapped/stored as
# plain text, matching every other row already in the sheet.
# Column C (Total Bets) is numeric for the "fixed" rows (437-461), except
# for the final row (462) where the reset-after-cashout bug resurfaces and
# the value again lands in the sheet as text.
$rows = @(
    @(437, "6.62", "2024-08-11 20:52:17", "5438", "n"),
    @(438, "5.15", "2024-08-11 21:00:36", "4149", "n"),
    @(439, "1.00", "2024-08-11 21:00:53", "4350", "n"),
    @(440, "2.57", "2024-08-11 21:01:22", "3883", "n"),
    @(441, "1.69", "2024-08-11 21:01:48", "4508", "n"),
    @(442, "4.93", "2024-08-11 21:05:43", "4344", "n"),
    @(443, "4.93", "2024-08-11 21:06:09", "4236", "n"),
    @(444, "4.96", "2024-08-11 21:11:19", "4797", "n"),
    @(445, "1.20", "2024-08-11 21:11:39", "4919", "n"),
    @(446, "1.70", "2024-08-11 21:12:04", "4150", "n"),
    @(447, "1.14", "2024-08-11 21:12:21", "4721", "n"),
    @(448, "1.18", "2024-08-11 21:14:16", "2926", "n"),
    @(449, "1.91", "2024-08-11 21:14:51", "3825", "n"),
    @(450, "1.04", "2024-08-11 21:17:50", "2767", "n"),
    @(451, "1.97", "2024-08-11 21:18:35", "3101", "n"),
    @(452, "1.11", "2024-08-11 21:18:53", "2910", "n"),
    @(453, "347.84", "2024-08-11 21:24:34", "3955", "n"),
    @(454, "2.44", "2024-08-11 21:25:05", "4049", "n"),
    @(455, "1.11", "2024-08-11 21:25:25", "4292", "n"),
    @(456, "1.11", "2024-08-11 21:25:26", "4292", "n"),
    @(457, "1.36", "2024-08-11 21:26:24", "3323", "n"),
    @(458, "1.64", "2024-08-11 21:26:59", "4589", "n"),
    @(459, "1.58", "2024-08-11 21:27:21", "4009", "n"),
    @(460, "1.58", "2024-08-11 21:27:57", "3791", "n"),
    @(461, "2.75", "2024-08-11 21:28:36", "4731", "n"),
    @(462, "1.13", "2024-08-11 21:28:53", "4690", "s")
)

foreach ($row in $rows) {
    $r = $row[0]
    $multiplier = $row[1]
    $timestamp = $row[2]
    $totalBets = $row[3]
    $totalBetsKind = $row[4]

    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $multiplier

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $timestamp

    if ($totalBetsKind -eq "n") {
        $ws.Cells.Item($r, 3).Value = [double]$totalBets
    } else {
        $ws.Cells.Item($r, 3).NumberFormat = "@"
        $ws.Cells.Item($r, 3).Value = $totalBets
    }
}
